$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.303.83'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '3.025.69'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'577.72"
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").Value = "'168.70"
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.024.54'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("D12").Value = "'0.477"
$ws.Range("E12").Value = '  +4.10%  '
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").Value = "'36.96"
$ws.Range("E14").Value = '  +6.09%  '
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = '66.365.42'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '3.510.78'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = "'7.24"
$ws.Range("E18").Value = '  +4.01%  '
$ws.Range("D19").Value = '3.015.37'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = "'16.26"
$ws.Range("E20").Value = '  +16.41%  '
$ws.Range("D21").Value = "'463.23"
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").Value = "'0.708"
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("D23").Value = "'7.50"
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").Value = "'83.53"
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = "'12.73"
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").Value = "'10.30"
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("D29").Value = "'8.44"
$ws.Range("E29").Value = '  +4.42%  '
$ws.Range("E30").Value = '  +4.72%  '
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").Value = "'0.0000102"
$ws.Range("E32").Value = '  -3.90%  '
$ws.Range("D33").Value = "'0.119"
$ws.Range("E33").Value = '  +6.58%  '
$ws.Range("D34").Value = "'28.12"
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = "'5.85"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = "'48.23"
$ws.Range("E38").Value = '  +11.23%  '
$ws.Range("D39").Value = "'2.06"
$ws.Range("E39").Value = '  -6.25%  '
$ws.Range("D40").Value = "'49.60"
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("D41").Value = "'0.313"
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = "'8.65"
$ws.Range("E44").Value = '  +2.08%  '
$ws.Range("D45").Value = "'386.34"
$ws.Range("E45").Value = '  -2.79%  '
$ws.Range("D46").Value = "'0.0359"
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '2.722.16'
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").Value = "'133.78"
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = "'24.79"
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").Value = "'2.25"
$ws.Range("E51").Value = '  +3.58%  '
